$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 17.6866637779875
$ws.Cells.Item(2, 3).Value = 10.6983053153646
$ws.Cells.Item(2, 4).Value = 13.83711915403442
$ws.Cells.Item(2, 5).Value = 14.65007535403013
$ws.Cells.Item(2, 7).Value = 39.18698738840938
$ws.Cells.Item(2, 8).Value = 16.99749112202236
$ws.Cells.Item(2, 10).Value = 8.910433531478912
$ws.Cells.Item(2, 12).Value = 11.83112848282943
$ws.Cells.Item(2, 13).Value = 17.15780910989969
$ws.Cells.Item(2, 15).Value = 27.24295665756311
$ws.Cells.Item(3, 2).Value = 17.22427569990756
$ws.Cells.Item(3, 3).Value = 10.5331417801198
$ws.Cells.Item(3, 4).Value = 13.86350227562167
$ws.Cells.Item(3, 5).Value = 14.6998951428765
$ws.Cells.Item(3, 7).Value = 39.34762155926313
$ws.Cells.Item(3, 8).Value = 17.06374967402619
$ws.Cells.Item(3, 10).Value = 8.916098445836479
$ws.Cells.Item(3, 12).Value = 11.82997651747835
$ws.Cells.Item(3, 13).Value = 17.04850186824684
$ws.Cells.Item(3, 15).Value = 27.35813233340279
$ws.Cells.Item(4, 2).Value = 16.93553387802879
$ws.Cells.Item(4, 3).Value = 10.43036066999414
$ws.Cells.Item(4, 4).Value = 13.88216668703177
$ws.Cells.Item(4, 5).Value = 14.73222460620112
$ws.Cells.Item(4, 7).Value = 39.45853227954999
$ws.Cells.Item(4, 8).Value = 17.10744027263861
$ws.Cells.Item(4, 10).Value = 8.919857241285207
$ws.Cells.Item(4, 12).Value = 11.83054290115346
$ws.Cells.Item(4, 13).Value = 16.98283087974408
$ws.Cells.Item(4, 15).Value = 27.43503725829338
$ws.Cells.Item(5, 2).Value = 16.81681645472039
$ws.Cells.Item(5, 3).Value = 10.38816969930818
$ws.Cells.Item(5, 4).Value = 13.89039217764852
$ws.Cells.Item(5, 5).Value = 14.74583762431645
$ws.Cells.Item(5, 7).Value = 39.50680323969567
$ws.Cells.Item(5, 8).Value = 17.12600066305139
$ws.Cells.Item(5, 10).Value = 8.921459743143533
$ws.Cells.Item(5, 12).Value = 11.83109514089269
$ws.Cells.Item(5, 13).Value = 16.95645151891427
$ws.Cells.Item(5, 15).Value = 27.467929246746
$ws.Cells.Item(6, 2).Value = 16.79704486625934
$ws.Cells.Item(6, 3).Value = 10.38114650663763
$ws.Cells.Item(6, 4).Value = 13.8917954261656
$ws.Cells.Item(6, 5).Value = 14.74812456892522
$ws.Cells.Item(6, 7).Value = 39.5150038312383
$ws.Cells.Item(6, 8).Value = 17.12912826482741
$ws.Cells.Item(6, 10).Value = 8.921730117823518
$ws.Cells.Item(6, 12).Value = 11.83120628585777
$ws.Cells.Item(6, 13).Value = 16.95209486748038
$ws.Cells.Item(6, 15).Value = 27.47348462006332
$ws.Cells.Item(7, 2).Value = 16.93393685429277
$ws.Cells.Item(7, 3).Value = 10.42979285975921
$ws.Cells.Item(7, 4).Value = 13.88227511059467
$ws.Cells.Item(7, 5).Value = 14.73240641934538
$ws.Cells.Item(7, 7).Value = 39.45917085142477
$ws.Cells.Item(7, 8).Value = 17.10768752328543
$ws.Cells.Item(7, 10).Value = 8.919878566375129
$ws.Cells.Item(7, 12).Value = 11.83054904591181
$ws.Cells.Item(7, 13).Value = 16.98247354670986
$ws.Cells.Item(7, 15).Value = 27.43547456879217
$ws.Cells.Item(8, 2).Value = 17.52833028010949
$ws.Cells.Item(8, 3).Value = 10.64166359261816
$ws.Cells.Item(8, 4).Value = 13.84570429354683
$ws.Cells.Item(8, 5).Value = 14.66689269184728
$ws.Cells.Item(8, 7).Value = 39.23981684574828
$ws.Cells.Item(8, 8).Value = 17.01971274892228
$ws.Cells.Item(8, 10).Value = 8.912328707460658
$ws.Cells.Item(8, 12).Value = 11.83046764169889
$ws.Cells.Item(8, 13).Value = 17.11983178774524
$ws.Cells.Item(8, 15).Value = 27.28138337265275
$ws.Cells.Item(9, 2).Value = 18.6489243222035
$ws.Cells.Item(9, 3).Value = 11.04466716589779
$ws.Cells.Item(9, 4).Value = 13.79355753219077
$ws.Cells.Item(9, 5).Value = 14.55218093149132
$ws.Cells.Item(9, 7).Value = 38.90771523411095
$ws.Cells.Item(9, 8).Value = 16.87106550964653
$ws.Cells.Item(9, 10).Value = 8.899739340492427
$ws.Cells.Item(9, 12).Value = 11.84036170245179
$ws.Cells.Item(9, 13).Value = 17.39979450224508
$ws.Cells.Item(9, 15).Value = 27.02843451194385
$ws.Cells.Item(10, 2).Value = 19.43659359980696
$ws.Cells.Item(10, 3).Value = 11.33108454314236
$ws.Cells.Item(10, 4).Value = 13.76718325948352
$ws.Cells.Item(10, 5).Value = 14.47622798226169
$ws.Cells.Item(10, 7).Value = 38.72429192338888
$ws.Cells.Item(10, 8).Value = 16.77641658569591
$ws.Cells.Item(10, 10).Value = 8.891827480289027
$ws.Cells.Item(10, 12).Value = 11.85368663155074
$ws.Cells.Item(10, 13).Value = 17.61077638388364
$ws.Cells.Item(10, 15).Value = 26.87278464110016
$ws.Cells.Item(11, 2).Value = 19.78566594983639
$ws.Cells.Item(11, 3).Value = 11.45887246487156
$ws.Cells.Item(11, 4).Value = 13.75777686697511
$ws.Cells.Item(11, 5).Value = 14.44346958857852
$ws.Cells.Item(11, 7).Value = 38.65415810737455
$ws.Cells.Item(11, 8).Value = 16.73652283605464
$ws.Cells.Item(11, 10).Value = 8.888515804695601
$ws.Cells.Item(11, 12).Value = 11.86104623908867
$ws.Cells.Item(11, 13).Value = 17.70765081692322
$ws.Cells.Item(11, 15).Value = 26.80857017327761
$ws.Cells.Item(12, 2).Value = 19.91640883589585
$ws.Cells.Item(12, 3).Value = 11.50686979804277
$ws.Cells.Item(12, 4).Value = 13.75458731378463
$ws.Cells.Item(12, 5).Value = 14.43132168384667
$ws.Cells.Item(12, 7).Value = 38.62952554191521
$ws.Cells.Item(12, 8).Value = 16.72187118220873
$ws.Cells.Item(12, 10).Value = 8.887302873645044
$ws.Cells.Item(12, 12).Value = 11.86401815699884
$ws.Cells.Item(12, 13).Value = 17.74444176754536
$ws.Cells.Item(12, 15).Value = 26.78520489110682
$ws.Cells.Item(13, 2).Value = 19.8883169134756
$ws.Cells.Item(13, 3).Value = 11.49655070858248
$ws.Cells.Item(13, 4).Value = 13.75525768102761
$ws.Cells.Item(13, 5).Value = 14.43392653613389
$ws.Cells.Item(13, 7).Value = 38.63474477959456
$ws.Cells.Item(13, 8).Value = 16.72500642320812
$ws.Cells.Item(13, 10).Value = 8.88756227415613
$ws.Cells.Item(13, 12).Value = 11.86336990100027
$ws.Cells.Item(13, 13).Value = 17.73651378979768
$ws.Cells.Item(13, 15).Value = 26.79019465620099
$ws.Cells.Item(14, 2).Value = 19.79645178791466
$ws.Cells.Item(14, 3).Value = 11.46282929754084
$ws.Cells.Item(14, 4).Value = 13.75750699774675
$ws.Cells.Item(14, 5).Value = 14.44246502749542
$ws.Cells.Item(14, 7).Value = 38.6520929129722
$ws.Cells.Item(14, 8).Value = 16.73530830914144
$ws.Cells.Item(14, 10).Value = 8.888415193041755
$ws.Cells.Item(14, 12).Value = 11.86128704282391
$ws.Cells.Item(14, 13).Value = 17.71067562062419
$ws.Cells.Item(14, 15).Value = 26.80662881123499
$ws.Cells.Item(15, 2).Value = 19.73999058053106
$ws.Cells.Item(15, 3).Value = 11.44212180415677
$ws.Cells.Item(15, 4).Value = 13.75893326443192
$ws.Cells.Item(15, 5).Value = 14.44772854267678
$ws.Cells.Item(15, 7).Value = 38.6629702610524
$ws.Cells.Item(15, 8).Value = 16.74167781223758
$ws.Cells.Item(15, 10).Value = 8.888942980466409
$ws.Cells.Item(15, 12).Value = 11.86003527299316
$ws.Cells.Item(15, 13).Value = 17.69486222594846
$ws.Cells.Item(15, 15).Value = 26.81681920177399
$ws.Cells.Item(16, 2).Value = 19.41358492287191
$ws.Cells.Item(16, 3).Value = 11.32268004240967
$ws.Cells.Item(16, 4).Value = 13.76785012087451
$ws.Cells.Item(16, 5).Value = 14.47840482758546
$ws.Cells.Item(16, 7).Value = 38.7291441478767
$ws.Cells.Item(16, 8).Value = 16.77908741449308
$ws.Cells.Item(16, 10).Value = 8.892049672785372
$ws.Cells.Item(16, 12).Value = 11.8532316435784
$ws.Cells.Item(16, 13).Value = 17.6044615653506
$ws.Cells.Item(16, 15).Value = 26.87711414469013
$ws.Cells.Item(17, 2).Value = 19.21089388933375
$ws.Cells.Item(17, 3).Value = 11.2487408759081
$ws.Cells.Item(17, 4).Value = 13.77398394104332
$ws.Cells.Item(17, 5).Value = 14.49768235873657
$ws.Cells.Item(17, 7).Value = 38.77315665674003
$ws.Cells.Item(17, 8).Value = 16.80284732600427
$ws.Cells.Item(17, 10).Value = 8.894029006425082
$ws.Cells.Item(17, 12).Value = 11.84938913981403
$ws.Cells.Item(17, 13).Value = 17.5492177284025
$ws.Cells.Item(17, 15).Value = 26.91579387205744
$ws.Cells.Item(18, 2).Value = 19.09344759314907
$ws.Cells.Item(18, 3).Value = 11.20597913545036
$ws.Cells.Item(18, 4).Value = 13.77775588600895
$ws.Cells.Item(18, 5).Value = 14.50893909840268
$ws.Cells.Item(18, 7).Value = 38.79972326773583
$ws.Cells.Item(18, 8).Value = 16.81681104064304
$ws.Cells.Item(18, 10).Value = 8.895194537156975
$ws.Cells.Item(18, 12).Value = 11.84730131492229
$ws.Cells.Item(18, 13).Value = 17.51752878142951
$ws.Cells.Item(18, 15).Value = 26.93866163902313
$ws.Cells.Item(19, 2).Value = 19.05353759051302
$ws.Cells.Item(19, 3).Value = 11.19146160590652
$ws.Cells.Item(19, 4).Value = 13.77907490163098
$ws.Cells.Item(19, 5).Value = 14.51277945676574
$ws.Cells.Item(19, 7).Value = 38.808932877185
$ws.Cells.Item(19, 8).Value = 16.82159002050444
$ws.Cells.Item(19, 10).Value = 8.895593821596556
$ws.Cells.Item(19, 12).Value = 11.84661546446616
$ws.Cells.Item(19, 13).Value = 17.50681487282339
$ws.Cells.Item(19, 15).Value = 26.94651068633424
$ws.Cells.Item(20, 2).Value = 19.23256095328922
$ws.Cells.Item(20, 3).Value = 11.25663627696303
$ws.Cells.Item(20, 4).Value = 13.77330574021802
$ws.Cells.Item(20, 5).Value = 14.49561276718848
$ws.Cells.Item(20, 7).Value = 38.76834180495635
$ws.Cells.Item(20, 8).Value = 16.80028723461181
$ws.Cells.Item(20, 10).Value = 8.893815502816336
$ws.Cells.Item(20, 12).Value = 11.84978553696581
$ws.Cells.Item(20, 13).Value = 17.55508980610541
$ws.Cells.Item(20, 15).Value = 26.9116121340054
$ws.Cells.Item(21, 2).Value = 19.82347483025382
$ws.Cells.Item(21, 3).Value = 11.47274501914705
$ws.Cells.Item(21, 4).Value = 13.75683621305824
$ws.Cells.Item(21, 5).Value = 14.4399500976661
$ws.Cells.Item(21, 7).Value = 38.64694499288098
$ws.Cells.Item(21, 8).Value = 16.73227003648043
$ws.Cells.Item(21, 10).Value = 8.888163555530113
$ws.Cells.Item(21, 12).Value = 11.86189382173396
$ws.Cells.Item(21, 13).Value = 17.71826220019194
$ws.Cells.Item(21, 15).Value = 26.80177585468689
$ws.Cells.Item(22, 2).Value = 20.20121312935593
$ws.Cells.Item(22, 3).Value = 11.61167866319146
$ws.Cells.Item(22, 4).Value = 13.74824313348332
$ws.Cells.Item(22, 5).Value = 14.4050688667175
$ws.Cells.Item(22, 7).Value = 38.57883402610739
$ws.Cells.Item(22, 8).Value = 16.69047060290613
$ws.Cells.Item(22, 10).Value = 8.88470931435511
$ws.Cells.Item(22, 12).Value = 11.87088474255038
$ws.Cells.Item(22, 13).Value = 17.82551631350238
$ws.Cells.Item(22, 15).Value = 26.73553824124554
$ws.Cells.Item(23, 2).Value = 20.00041609059504
$ws.Cells.Item(23, 3).Value = 11.53774863679593
$ws.Cells.Item(23, 4).Value = 13.75263089810343
$ws.Cells.Item(23, 5).Value = 14.42354888563677
$ws.Cells.Item(23, 7).Value = 38.61415488317779
$ws.Cells.Item(23, 8).Value = 16.71253676912125
$ws.Cells.Item(23, 10).Value = 8.886531050144256
$ws.Cells.Item(23, 12).Value = 11.86598810255159
$ws.Cells.Item(23, 13).Value = 17.7682242384418
$ws.Cells.Item(23, 15).Value = 26.77038183073068
$ws.Cells.Item(24, 2).Value = 19.22276811964994
$ws.Cells.Item(24, 3).Value = 11.25306755208038
$ws.Cells.Item(24, 4).Value = 13.77361159000162
$ws.Cells.Item(24, 5).Value = 14.49654788799375
$ws.Cells.Item(24, 7).Value = 38.77051466544883
$ws.Cells.Item(24, 8).Value = 16.80144370554861
$ws.Cells.Item(24, 10).Value = 8.893911941853446
$ws.Cells.Item(24, 12).Value = 11.84960594791009
$ws.Cells.Item(24, 13).Value = 17.55243481544099
$ws.Cells.Item(24, 15).Value = 26.91350073474016
$ws.Cells.Item(25, 2).Value = 18.35147923310608
$ws.Cells.Item(25, 3).Value = 10.93719760878318
$ws.Cells.Item(25, 4).Value = 13.80556836730898
$ws.Cells.Item(25, 5).Value = 14.58174690327507
$ws.Cells.Item(25, 7).Value = 38.98697570637674
$ws.Cells.Item(25, 8).Value = 16.90872201743242
$ws.Cells.Item(25, 10).Value = 8.902909254139328
$ws.Cells.Item(25, 12).Value = 11.83661600371732
$ws.Cells.Item(25, 13).Value = 17.32304136923039
$ws.Cells.Item(25, 15).Value = 27.09157507568348
